$d = $word.ActiveDocument

# Update the date paragraph (first paragraph, outside the table)
$dateRange = $d.Paragraphs.Item(1).Range
$dateRange.Find.Execute("2025-08-30 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-31 Sunday", 1) | Out-Null

# Update each math expression cell in the table, addressed by (row, column)
# so that duplicate expressions in different cells get the correct, distinct replacement.
$t = $d.Tables.Item(1)

$cellRange = $t.Cell(1, 1).Range
$cellRange.Find.Execute("9+73=82", $true, $false, $false, $false, $false, $true, 1, $false, "87-44=43", 1) | Out-Null
$cellRange = $t.Cell(1, 2).Range
$cellRange.Find.Execute("65-15=50", $true, $false, $false, $false, $false, $true, 1, $false, "88-74=14", 1) | Out-Null
$cellRange = $t.Cell(1, 3).Range
$cellRange.Find.Execute("59-28=31", $true, $false, $false, $false, $false, $true, 1, $false, "59+3=62", 1) | Out-Null
$cellRange = $t.Cell(1, 4).Range
$cellRange.Find.Execute("71-68=3", $true, $false, $false, $false, $false, $true, 1, $false, "24+10=34", 1) | Out-Null
$cellRange = $t.Cell(1, 5).Range
$cellRange.Find.Execute("82-68=14", $true, $false, $false, $false, $false, $true, 1, $false, "85-51=34", 1) | Out-Null
$cellRange = $t.Cell(2, 1).Range
$cellRange.Find.Execute("70-21=49", $true, $false, $false, $false, $false, $true, 1, $false, "60-50=10", 1) | Out-Null
$cellRange = $t.Cell(2, 2).Range
$cellRange.Find.Execute("0+19=19", $true, $false, $false, $false, $false, $true, 1, $false, "84-28=56", 1) | Out-Null
$cellRange = $t.Cell(2, 3).Range
$cellRange.Find.Execute("88-65=23", $true, $false, $false, $false, $false, $true, 1, $false, "73-6=67", 1) | Out-Null
$cellRange = $t.Cell(2, 4).Range
$cellRange.Find.Execute("81-78=3", $true, $false, $false, $false, $false, $true, 1, $false, "67-5=62", 1) | Out-Null
$cellRange = $t.Cell(2, 5).Range
$cellRange.Find.Execute("10+3=13", $true, $false, $false, $false, $false, $true, 1, $false, "28+68=96", 1) | Out-Null
$cellRange = $t.Cell(3, 1).Range
$cellRange.Find.Execute("59+18=77", $true, $false, $false, $false, $false, $true, 1, $false, "48-24=24", 1) | Out-Null
$cellRange = $t.Cell(3, 2).Range
$cellRange.Find.Execute("46+15=61", $true, $false, $false, $false, $false, $true, 1, $false, "95+2=97", 1) | Out-Null
$cellRange = $t.Cell(3, 3).Range
$cellRange.Find.Execute("20+5=25", $true, $false, $false, $false, $false, $true, 1, $false, "18+49=67", 1) | Out-Null
$cellRange = $t.Cell(3, 4).Range
$cellRange.Find.Execute("54-16=38", $true, $false, $false, $false, $false, $true, 1, $false, "3+24=27", 1) | Out-Null
$cellRange = $t.Cell(3, 5).Range
$cellRange.Find.Execute("98-45=53", $true, $false, $false, $false, $false, $true, 1, $false, "36-5=31", 1) | Out-Null
$cellRange = $t.Cell(4, 1).Range
$cellRange.Find.Execute("34-23=11", $true, $false, $false, $false, $false, $true, 1, $false, "89+8=97", 1) | Out-Null
$cellRange = $t.Cell(4, 2).Range
$cellRange.Find.Execute("44+8=52", $true, $false, $false, $false, $false, $true, 1, $false, "45-17=28", 1) | Out-Null
$cellRange = $t.Cell(4, 3).Range
$cellRange.Find.Execute("39+46=85", $true, $false, $false, $false, $false, $true, 1, $false, "33-3=30", 1) | Out-Null
$cellRange = $t.Cell(4, 4).Range
$cellRange.Find.Execute("2+31=33", $true, $false, $false, $false, $false, $true, 1, $false, "57-1=56", 1) | Out-Null
$cellRange = $t.Cell(4, 5).Range
$cellRange.Find.Execute("6-3=3", $true, $false, $false, $false, $false, $true, 1, $false, "16+62=78", 1) | Out-Null
$cellRange = $t.Cell(5, 1).Range
$cellRange.Find.Execute("42-10=32", $true, $false, $false, $false, $false, $true, 1, $false, "47-20=27", 1) | Out-Null
$cellRange = $t.Cell(5, 2).Range
$cellRange.Find.Execute("3+72=75", $true, $false, $false, $false, $false, $true, 1, $false, "43-25=18", 1) | Out-Null
$cellRange = $t.Cell(5, 3).Range
$cellRange.Find.Execute("2+3=5", $true, $false, $false, $false, $false, $true, 1, $false, "98-54=44", 1) | Out-Null
$cellRange = $t.Cell(5, 4).Range
$cellRange.Find.Execute("14+23=37", $true, $false, $false, $false, $false, $true, 1, $false, "67-30=37", 1) | Out-Null
$cellRange = $t.Cell(5, 5).Range
$cellRange.Find.Execute("25+55=80", $true, $false, $false, $false, $false, $true, 1, $false, "59-31=28", 1) | Out-Null
$cellRange = $t.Cell(6, 1).Range
$cellRange.Find.Execute("59-46=13", $true, $false, $false, $false, $false, $true, 1, $false, "36+47=83", 1) | Out-Null
$cellRange = $t.Cell(6, 2).Range
$cellRange.Find.Execute("99-37=62", $true, $false, $false, $false, $false, $true, 1, $false, "87-35=52", 1) | Out-Null
$cellRange = $t.Cell(6, 3).Range
$cellRange.Find.Execute("16+38=54", $true, $false, $false, $false, $false, $true, 1, $false, "74-59=15", 1) | Out-Null
$cellRange = $t.Cell(6, 4).Range
$cellRange.Find.Execute("78-7=71", $true, $false, $false, $false, $false, $true, 1, $false, "56-7=49", 1) | Out-Null
$cellRange = $t.Cell(6, 5).Range
$cellRange.Find.Execute("54-33=21", $true, $false, $false, $false, $false, $true, 1, $false, "18+30=48", 1) | Out-Null
$cellRange = $t.Cell(7, 1).Range
$cellRange.Find.Execute("34-19=15", $true, $false, $false, $false, $false, $true, 1, $false, "64+13=77", 1) | Out-Null
$cellRange = $t.Cell(7, 2).Range
$cellRange.Find.Execute("45+22=67", $true, $false, $false, $false, $false, $true, 1, $false, "89-62=27", 1) | Out-Null
$cellRange = $t.Cell(7, 3).Range
$cellRange.Find.Execute("83+8=91", $true, $false, $false, $false, $false, $true, 1, $false, "36-19=17", 1) | Out-Null
$cellRange = $t.Cell(7, 4).Range
$cellRange.Find.Execute("23-15=8", $true, $false, $false, $false, $false, $true, 1, $false, "78-51=27", 1) | Out-Null
$cellRange = $t.Cell(7, 5).Range
$cellRange.Find.Execute("37+29=66", $true, $false, $false, $false, $false, $true, 1, $false, "27+27=54", 1) | Out-Null
$cellRange = $t.Cell(8, 1).Range
$cellRange.Find.Execute("45+3=48", $true, $false, $false, $false, $false, $true, 1, $false, "65-1=64", 1) | Out-Null
$cellRange = $t.Cell(8, 2).Range
$cellRange.Find.Execute("93-89=4", $true, $false, $false, $false, $false, $true, 1, $false, "97-9=88", 1) | Out-Null
$cellRange = $t.Cell(8, 3).Range
$cellRange.Find.Execute("86-82=4", $true, $false, $false, $false, $false, $true, 1, $false, "39-31=8", 1) | Out-Null
$cellRange = $t.Cell(8, 4).Range
$cellRange.Find.Execute("35+30=65", $true, $false, $false, $false, $false, $true, 1, $false, "17+20=37", 1) | Out-Null
$cellRange = $t.Cell(8, 5).Range
$cellRange.Find.Execute("73-47=26", $true, $false, $false, $false, $false, $true, 1, $false, "22-9=13", 1) | Out-Null
$cellRange = $t.Cell(9, 1).Range
$cellRange.Find.Execute("67+3=70", $true, $false, $false, $false, $false, $true, 1, $false, "96-31=65", 1) | Out-Null
$cellRange = $t.Cell(9, 2).Range
$cellRange.Find.Execute("55-29=26", $true, $false, $false, $false, $false, $true, 1, $false, "33-3=30", 1) | Out-Null
$cellRange = $t.Cell(9, 3).Range
$cellRange.Find.Execute("41-9=32", $true, $false, $false, $false, $false, $true, 1, $false, "71-11=60", 1) | Out-Null
$cellRange = $t.Cell(9, 4).Range
$cellRange.Find.Execute("2+27=29", $true, $false, $false, $false, $false, $true, 1, $false, "61+34=95", 1) | Out-Null
$cellRange = $t.Cell(9, 5).Range
$cellRange.Find.Execute("8+87=95", $true, $false, $false, $false, $false, $true, 1, $false, "33-6=27", 1) | Out-Null
$cellRange = $t.Cell(10, 1).Range
$cellRange.Find.Execute("52+13=65", $true, $false, $false, $false, $false, $true, 1, $false, "93-77=16", 1) | Out-Null
$cellRange = $t.Cell(10, 2).Range
$cellRange.Find.Execute("81+5=86", $true, $false, $false, $false, $false, $true, 1, $false, "67-1=66", 1) | Out-Null
$cellRange = $t.Cell(10, 3).Range
$cellRange.Find.Execute("42+8=50", $true, $false, $false, $false, $false, $true, 1, $false, "51-8=43", 1) | Out-Null
$cellRange = $t.Cell(10, 4).Range
$cellRange.Find.Execute("45-9=36", $true, $false, $false, $false, $false, $true, 1, $false, "84-9=75", 1) | Out-Null
$cellRange = $t.Cell(10, 5).Range
$cellRange.Find.Execute("12+23=35", $true, $false, $false, $false, $false, $true, 1, $false, "66-13=53", 1) | Out-Null
$cellRange = $t.Cell(11, 1).Range
$cellRange.Find.Execute("71-62=9", $true, $false, $false, $false, $false, $true, 1, $false, "91-82=9", 1) | Out-Null
$cellRange = $t.Cell(11, 2).Range
$cellRange.Find.Execute("71+25=96", $true, $false, $false, $false, $false, $true, 1, $false, "3+64=67", 1) | Out-Null
$cellRange = $t.Cell(11, 3).Range
$cellRange.Find.Execute("80-64=16", $true, $false, $false, $false, $false, $true, 1, $false, "91-89=2", 1) | Out-Null
$cellRange = $t.Cell(11, 4).Range
$cellRange.Find.Execute("24+39=63", $true, $false, $false, $false, $false, $true, 1, $false, "72+19=91", 1) | Out-Null
$cellRange = $t.Cell(11, 5).Range
$cellRange.Find.Execute("26-20=6", $true, $false, $false, $false, $false, $true, 1, $false, "21+74=95", 1) | Out-Null
$cellRange = $t.Cell(12, 1).Range
$cellRange.Find.Execute("54-9=45", $true, $false, $false, $false, $false, $true, 1, $false, "56-20=36", 1) | Out-Null
$cellRange = $t.Cell(12, 2).Range
$cellRange.Find.Execute("72+8=80", $true, $false, $false, $false, $false, $true, 1, $false, "57-9=48", 1) | Out-Null
$cellRange = $t.Cell(12, 3).Range
$cellRange.Find.Execute("82-11=71", $true, $false, $false, $false, $false, $true, 1, $false, "74-58=16", 1) | Out-Null
$cellRange = $t.Cell(12, 4).Range
$cellRange.Find.Execute("52-7=45", $true, $false, $false, $false, $false, $true, 1, $false, "79+19=98", 1) | Out-Null
$cellRange = $t.Cell(12, 5).Range
$cellRange.Find.Execute("83-2=81", $true, $false, $false, $false, $false, $true, 1, $false, "85-82=3", 1) | Out-Null
$cellRange = $t.Cell(13, 1).Range
$cellRange.Find.Execute("55+43=98", $true, $false, $false, $false, $false, $true, 1, $false, "11+46=57", 1) | Out-Null
$cellRange = $t.Cell(13, 2).Range
$cellRange.Find.Execute("38-18=20", $true, $false, $false, $false, $false, $true, 1, $false, "33+59=92", 1) | Out-Null
$cellRange = $t.Cell(13, 3).Range
$cellRange.Find.Execute("33-1=32", $true, $false, $false, $false, $false, $true, 1, $false, "23+60=83", 1) | Out-Null
$cellRange = $t.Cell(13, 4).Range
$cellRange.Find.Execute("54+39=93", $true, $false, $false, $false, $false, $true, 1, $false, "92-14=78", 1) | Out-Null
$cellRange = $t.Cell(13, 5).Range
$cellRange.Find.Execute("90-29=61", $true, $false, $false, $false, $false, $true, 1, $false, "76+19=95", 1) | Out-Null
$cellRange = $t.Cell(14, 1).Range
$cellRange.Find.Execute("34-25=9", $true, $false, $false, $false, $false, $true, 1, $false, "30+7=37", 1) | Out-Null
$cellRange = $t.Cell(14, 2).Range
$cellRange.Find.Execute("56+17=73", $true, $false, $false, $false, $false, $true, 1, $false, "32-11=21", 1) | Out-Null
$cellRange = $t.Cell(14, 3).Range
$cellRange.Find.Execute("31+50=81", $true, $false, $false, $false, $false, $true, 1, $false, "25-17=8", 1) | Out-Null
$cellRange = $t.Cell(14, 4).Range
$cellRange.Find.Execute("90+3=93", $true, $false, $false, $false, $false, $true, 1, $false, "92-78=14", 1) | Out-Null
$cellRange = $t.Cell(14, 5).Range
$cellRange.Find.Execute("38-23=15", $true, $false, $false, $false, $false, $true, 1, $false, "41-2=39", 1) | Out-Null
$cellRange = $t.Cell(15, 1).Range
$cellRange.Find.Execute("23-20=3", $true, $false, $false, $false, $false, $true, 1, $false, "88-83=5", 1) | Out-Null
$cellRange = $t.Cell(15, 2).Range
$cellRange.Find.Execute("65-15=50", $true, $false, $false, $false, $false, $true, 1, $false, "59-35=24", 1) | Out-Null
$cellRange = $t.Cell(15, 3).Range
$cellRange.Find.Execute("60-28=32", $true, $false, $false, $false, $false, $true, 1, $false, "22+41=63", 1) | Out-Null
$cellRange = $t.Cell(15, 4).Range
$cellRange.Find.Execute("88+10=98", $true, $false, $false, $false, $false, $true, 1, $false, "90-36=54", 1) | Out-Null
$cellRange = $t.Cell(15, 5).Range
$cellRange.Find.Execute("90-31=59", $true, $false, $false, $false, $false, $true, 1, $false, "21+52=73", 1) | Out-Null
$cellRange = $t.Cell(16, 1).Range
$cellRange.Find.Execute("74-44=30", $true, $false, $false, $false, $false, $true, 1, $false, "68+6=74", 1) | Out-Null
$cellRange = $t.Cell(16, 2).Range
$cellRange.Find.Execute("96-35=61", $true, $false, $false, $false, $false, $true, 1, $false, "55+27=82", 1) | Out-Null
$cellRange = $t.Cell(16, 3).Range
$cellRange.Find.Execute("33+50=83", $true, $false, $false, $false, $false, $true, 1, $false, "95-66=29", 1) | Out-Null
$cellRange = $t.Cell(16, 4).Range
$cellRange.Find.Execute("52-1=51", $true, $false, $false, $false, $false, $true, 1, $false, "81-1=80", 1) | Out-Null
$cellRange = $t.Cell(16, 5).Range
$cellRange.Find.Execute("65+27=92", $true, $false, $false, $false, $false, $true, 1, $false, "37+0=37", 1) | Out-Null
$cellRange = $t.Cell(17, 1).Range
$cellRange.Find.Execute("93+3=96", $true, $false, $false, $false, $false, $true, 1, $false, "35+23=58", 1) | Out-Null
$cellRange = $t.Cell(17, 2).Range
$cellRange.Find.Execute("0+15=15", $true, $false, $false, $false, $false, $true, 1, $false, "3+74=77", 1) | Out-Null
$cellRange = $t.Cell(17, 3).Range
$cellRange.Find.Execute("13+14=27", $true, $false, $false, $false, $false, $true, 1, $false, "74+9=83", 1) | Out-Null
$cellRange = $t.Cell(17, 4).Range
$cellRange.Find.Execute("40-16=24", $true, $false, $false, $false, $false, $true, 1, $false, "92-27=65", 1) | Out-Null
$cellRange = $t.Cell(17, 5).Range
$cellRange.Find.Execute("63-2=61", $true, $false, $false, $false, $false, $true, 1, $false, "20-13=7", 1) | Out-Null
$cellRange = $t.Cell(18, 1).Range
$cellRange.Find.Execute("12+21=33", $true, $false, $false, $false, $false, $true, 1, $false, "25+62=87", 1) | Out-Null
$cellRange = $t.Cell(18, 2).Range
$cellRange.Find.Execute("48+21=69", $true, $false, $false, $false, $false, $true, 1, $false, "30+44=74", 1) | Out-Null
$cellRange = $t.Cell(18, 3).Range
$cellRange.Find.Execute("35-15=20", $true, $false, $false, $false, $false, $true, 1, $false, "11+69=80", 1) | Out-Null
$cellRange = $t.Cell(18, 4).Range
$cellRange.Find.Execute("11+51=62", $true, $false, $false, $false, $false, $true, 1, $false, "7+18=25", 1) | Out-Null
$cellRange = $t.Cell(18, 5).Range
$cellRange.Find.Execute("83-55=28", $true, $false, $false, $false, $false, $true, 1, $false, "37+58=95", 1) | Out-Null
$cellRange = $t.Cell(19, 1).Range
$cellRange.Find.Execute("37+32=69", $true, $false, $false, $false, $false, $true, 1, $false, "63-37=26", 1) | Out-Null
$cellRange = $t.Cell(19, 2).Range
$cellRange.Find.Execute("63-44=19", $true, $false, $false, $false, $false, $true, 1, $false, "79-27=52", 1) | Out-Null
$cellRange = $t.Cell(19, 3).Range
$cellRange.Find.Execute("12+53=65", $true, $false, $false, $false, $false, $true, 1, $false, "77+3=80", 1) | Out-Null
$cellRange = $t.Cell(19, 4).Range
$cellRange.Find.Execute("18+55=73", $true, $false, $false, $false, $false, $true, 1, $false, "20+12=32", 1) | Out-Null
$cellRange = $t.Cell(19, 5).Range
$cellRange.Find.Execute("59-51=8", $true, $false, $false, $false, $false, $true, 1, $false, "44-37=7", 1) | Out-Null
$cellRange = $t.Cell(20, 1).Range
$cellRange.Find.Execute("17+80=97", $true, $false, $false, $false, $false, $true, 1, $false, "6+44=50", 1) | Out-Null
$cellRange = $t.Cell(20, 2).Range
$cellRange.Find.Execute("60+37=97", $true, $false, $false, $false, $false, $true, 1, $false, "70-16=54", 1) | Out-Null
$cellRange = $t.Cell(20, 3).Range
$cellRange.Find.Execute("63+35=98", $true, $false, $false, $false, $false, $true, 1, $false, "82-31=51", 1) | Out-Null
$cellRange = $t.Cell(20, 4).Range
$cellRange.Find.Execute("92+5=97", $true, $false, $false, $false, $false, $true, 1, $false, "10+30=40", 1) | Out-Null
$cellRange = $t.Cell(20, 5).Range
$cellRange.Find.Execute("25+2=27", $true, $false, $false, $false, $false, $true, 1, $false, "78-52=26", 1) | Out-Null

Write-Output "Done"
